$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 55: A Real Smooth Move
$ws.Range("H55").Value = 221.46153
$ws.Range("I55").Value = 65.40000000000001
$ws.Range("J55").Value = 319
$ws.Range("K55").Value = 65.40000000000001
$ws.Range("L55").Value = 319
$ws.Range("M55").Value = 148.6
$ws.Range("N55").Value = -747

# Row 58: A Matter of Vital Importance
$ws.Range("H58").Value = 2636.2
$ws.Range("J58").Value = 6474.25
$ws.Range("L58").Value = 19422.75
$ws.Range("N58").Value = -19722.75

# Row 62: The Mustache Suits Him
$ws.Range("H62").Value = 3525.3845
$ws.Range("I62").Value = 3306.7778
$ws.Range("J62").Value = 4017.25
$ws.Range("K62").Value = 3306.7778
$ws.Range("L62").Value = 4017.25
$ws.Range("M62").Value = -2682.7778
$ws.Range("N62").Value = -5265.25

# Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 3525.3845
$ws.Range("I65").Value = 3306.7778
$ws.Range("J65").Value = 4017.25
$ws.Range("K65").Value = 16533.889
$ws.Range("L65").Value = 20086.25
$ws.Range("M65").Value = -13413.889
$ws.Range("N65").Value = -26326.25

# Row 116: Growing Up
$ws.Range("H116").Value = 9299.799999999999
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

# Row 127: Liquid Competence
$ws.Range("H127").Value = 12183.091
$ws.Range("I127").Value = 1127.375
$ws.Range("K127").Value = 3382.125
$ws.Range("M127").Value = 1577.875

# Row 136: I Like Big Brush and I Cannot Lie
$ws.Range("H136").Value = 77900
$ws.Range("J136").Value = 77900
$ws.Range("L136").Value = 77900
$ws.Range("N136").Value = -88100

# Row 138: All-night Crafting
$ws.Range("H138").Value = 2792.984
$ws.Range("I138").Value = 1377.619
$ws.Range("J138").Value = 3517.9268
$ws.Range("K138").Value = 4132.857
$ws.Range("L138").Value = 10553.7804
$ws.Range("M138").Value = 1007.143
$ws.Range("N138").Value = -20833.7804

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 15562.652
$ws.Range("I2").Value = 22603.268
$ws.Range("J2").Value = 2361.5
$ws.Range("K2").Value = 22603.268
$ws.Range("L2").Value = 2361.5
$ws.Range("M2").Value = -22490.268
$ws.Range("N2").Value = -2587.5

# Row 27: Get Me the Hard Stuff
$ws.Range("H27").Value = 101249.5
$ws.Range("J27").Value = 101249.5
$ws.Range("L27").Value = 101249.5
$ws.Range("N27").Value = -101617.5

# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 2693.5
$ws.Range("I45").Value = 1604
$ws.Range("J45").Value = 3783
$ws.Range("K45").Value = 1604
$ws.Range("L45").Value = 3783
$ws.Range("M45").Value = -1227
$ws.Range("N45").Value = -4537

# Row 116: No Scope
$ws.Range("H116").Value = 15562.652
$ws.Range("I116").Value = 22603.268
$ws.Range("J116").Value = 2361.5
$ws.Range("K116").Value = 22603.268
$ws.Range("L116").Value = 2361.5
$ws.Range("M116").Value = -20309.268
$ws.Range("N116").Value = -6949.5

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value = 15562.652
$ws.Range("I3").Value = 22603.268
$ws.Range("J3").Value = 2361.5
$ws.Range("K3").Value = 22603.268
$ws.Range("L3").Value = 2361.5
$ws.Range("M3").Value = -22489.268
$ws.Range("N3").Value = -2589.5

# Row 63: Punching Your Way to Success
$ws.Range("H63").Value = 50000
$ws.Range("J63").Value = 50000
$ws.Range("L63").Value = 50000
$ws.Range("N63").Value = -51372

# Row 66: Foreign Exchange (L)
$ws.Range("H66").Value = 50000
$ws.Range("J66").Value = 50000
$ws.Range("L66").Value = 150000
$ws.Range("N66").Value = -156864

# Row 102: Renting Mortality
$ws.Range("H102").Value = 17625
$ws.Range("I102").Value = 17625
$ws.Range("K102").Value = 17625
$ws.Range("M102").Value = -14380

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 3618.0557
$ws.Range("I105").Value = 1819.2222
$ws.Range("K105").Value = 1819.2222
$ws.Range("M105").Value = -72.22219999999993

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 1920.2593
$ws.Range("I58").Value = 1255.9474
$ws.Range("K58").Value = 1255.9474
$ws.Range("M58").Value = -1052.9474

# Row 99: O Pine
$ws.Range("H99").Value = 4057.2593
$ws.Range("I99").Value = 3766.5833
$ws.Range("J99").Value = 4289.8
$ws.Range("K99").Value = 3766.5833
$ws.Range("L99").Value = 4289.8
$ws.Range("M99").Value = -2268.5833
$ws.Range("N99").Value = -7285.8

# Row 126: A Better Conductor
$ws.Range("H126").Value = 4057.2593
$ws.Range("I126").Value = 3766.5833
$ws.Range("J126").Value = 4289.8
$ws.Range("K126").Value = 11299.7499
$ws.Range("L126").Value = 12869.4
$ws.Range("M126").Value = -8829.749899999999
$ws.Range("N126").Value = -17809.4

# Row 136: Turali Quality
$ws.Range("H136").Value = 1920.2593
$ws.Range("I136").Value = 1255.9474
$ws.Range("K136").Value = 3767.8422
$ws.Range("M136").Value = -1217.8422

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap
$ws.Range("H5").Value = 1158.5862
$ws.Range("I5").Value = 1162.5454
$ws.Range("K5").Value = 3487.6362
$ws.Range("M5").Value = -3375.6362

# Row 6: Meat-lover's Special
$ws.Range("H6").Value = 366.22223
$ws.Range("I6").Value = 349.2
$ws.Range("J6").Value = 387.5
$ws.Range("K6").Value = 1047.6
$ws.Range("L6").Value = 1162.5
$ws.Range("M6").Value = -934.5999999999999
$ws.Range("N6").Value = -1388.5

# Row 7: It's Always Sunny in Vylbrand
$ws.Range("H7").Value = 7248.2666
$ws.Range("I7").Value = 14542.857
$ws.Range("K7").Value = 43628.571
$ws.Range("M7").Value = -43516.571

# Row 22: A Total Nut Job
$ws.Range("H22").Value = 1122.6
$ws.Range("I22").Value = 853.25
$ws.Range("J22").Value = 2200
$ws.Range("K22").Value = 2559.75
$ws.Range("L22").Value = 6600
$ws.Range("M22").Value = -2390.75
$ws.Range("N22").Value = -6938

# Row 27: Brain Food
$ws.Range("H27").Value = 1122.6
$ws.Range("I27").Value = 853.25
$ws.Range("J27").Value = 2200
$ws.Range("K27").Value = 2559.75
$ws.Range("L27").Value = 6600
$ws.Range("M27").Value = -2457.75
$ws.Range("N27").Value = -6804

# Row 33: Cooking with Gas
$ws.Range("H33").Value = 30
$ws.Range("I33").Value = 30
$ws.Range("J33").Value = 30
$ws.Range("K33").Value = 180
$ws.Range("L33").Value = 180
$ws.Range("M33").Value = 103
$ws.Range("N33").Value = -746

# Row 39: Bloody Good Tart, This
$ws.Range("H39").Value = 2608
$ws.Range("J39").Value = 2608
$ws.Range("L39").Value = 7824
$ws.Range("N39").Value = -8412

# Row 56: Culture Club
$ws.Range("H56").Value = 7249.5557
$ws.Range("I56").Value = 7249.5557
$ws.Range("K56").Value = 7249.5557
$ws.Range("M56").Value = -6719.5557

# Row 98: Sweet Kiss of Death
$ws.Range("H98").Value = 2600
$ws.Range("I98").Value = 2600
$ws.Range("K98").Value = 7800
$ws.Range("M98").Value = -6302

# Row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 1158.5862
$ws.Range("I135").Value = 1162.5454
$ws.Range("K135").Value = 10462.9086
$ws.Range("M135").Value = -7927.908599999999

$ws = $wb.Worksheets.Item("GSM")
# Row 109: You're My Wonderhall
$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52080

# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 2566.3
$ws.Range("I122").Value = 2353.4
$ws.Range("J122").Value = 2992.1
$ws.Range("K122").Value = 7060.200000000001
$ws.Range("L122").Value = 8976.299999999999
$ws.Range("M122").Value = -4610.200000000001
$ws.Range("N122").Value = -13876.3

# Row 132: On Board for Lar
$ws.Range("H132").Value = 6176.846
$ws.Range("I132").Value = 6399.9165
$ws.Range("K132").Value = 19199.7495
$ws.Range("M132").Value = -16669.7495

# Row 134: Guaranteed Gem
$ws.Range("H134").Value = 76665.2
$ws.Range("J134").Value = 76665.2
$ws.Range("L134").Value = 229995.6
$ws.Range("N134").Value = -235065.6

$ws = $wb.Worksheets.Item("LTW")
# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 3624.5908
$ws.Range("I132").Value = 2337.7058
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 7013.117400000001
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -4483.117400000001
$ws.Range("N132").Value = -29060

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 2357.6
$ws.Range("I132").Value = 1697.0834
$ws.Range("K132").Value = 5091.2502
$ws.Range("M132").Value = -2561.2502
